$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price (D) and volume-change (E) values.
# Column D values are set via a text-formatted range to avoid Excel
# auto-converting numeric-looking strings (e.g. "578.41", "1.00") into
# floating point numbers, which would lose exact text formatting/precision.

$cell = $ws.Range('D2')
$cell.NumberFormat = "@"
$cell.Value = '61.768.83'
$cell.ClearFormats()
$ws.Range('E2').Value = '  -1.11%  '
$cell = $ws.Range('D3')
$cell.NumberFormat = "@"
$cell.Value = '2.445.36'
$cell.ClearFormats()
$ws.Range('E3').Value = '  -0.07%  '
$ws.Range('E4').Value = '  -0.11%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = "@"
$cell.Value = '578.41'
$cell.ClearFormats()
$ws.Range('E5').Value = '  -0.65%  '
$cell = $ws.Range('D6')
$cell.NumberFormat = "@"
$cell.Value = '140.73'
$cell.ClearFormats()
$ws.Range('E6').Value = '  -1.93%  '
$ws.Range('E7').Value = '  +0.09%  '
$cell = $ws.Range('D8')
$cell.NumberFormat = "@"
$cell.Value = '0.533'
$cell.ClearFormats()
$ws.Range('E8').Value = '  +0.75%  '
$cell = $ws.Range('D9')
$cell.NumberFormat = "@"
$cell.Value = '2.436.38'
$cell.ClearFormats()
$ws.Range('E9').Value = '  -0.36%  '
$ws.Range('E10').Value = '  +2.21%  '
$ws.Range('E11').Value = '  +1.97%  '
$ws.Range('E12').Value = '  -1.07%  '
$cell = $ws.Range('D13')
$cell.NumberFormat = "@"
$cell.Value = '0.340'
$cell.ClearFormats()
$ws.Range('E13').Value = '  -1.96%  '
$cell = $ws.Range('D14')
$cell.NumberFormat = "@"
$cell.Value = '25.89'
$cell.ClearFormats()
$ws.Range('E14').Value = '  -2.19%  '
$cell = $ws.Range('D15')
$cell.NumberFormat = "@"
$cell.Value = '2.906.35'
$cell.ClearFormats()
$ws.Range('E15').Value = '  +1.62%  '
$ws.Range('E16').Value = '  -1.00%  '
$cell = $ws.Range('D17')
$cell.NumberFormat = "@"
$cell.Value = '61.716.59'
$cell.ClearFormats()
$ws.Range('E17').Value = '  -0.89%  '
$cell = $ws.Range('D18')
$cell.NumberFormat = "@"
$cell.Value = '2.448.79'
$cell.ClearFormats()
$ws.Range('E18').Value = '  +0.57%  '
$ws.Range('E19').Value = '  -3.62%  '
$cell = $ws.Range('D20')
$cell.NumberFormat = "@"
$cell.Value = '7.23'
$cell.ClearFormats()
$ws.Range('E20').Value = '  +1.41%  '
$cell = $ws.Range('D21')
$cell.NumberFormat = "@"
$cell.Value = '325.06'
$cell.ClearFormats()
$ws.Range('E21').Value = '  -2.25%  '
$ws.Range('E22').Value = '  -1.17%  '
$cell = $ws.Range('D23')
$cell.NumberFormat = "@"
$cell.Value = '6.07'
$cell.ClearFormats()
$ws.Range('E23').Value = '  +1.27%  '
$ws.Range('E24').Value = '  -0.84%  '
$ws.Range('E25').Value = '  -0.08%  '
$cell = $ws.Range('D26')
$cell.NumberFormat = "@"
$cell.Value = '64.92'
$cell.ClearFormats()
$ws.Range('E26').Value = '  -1.53%  '
$cell = $ws.Range('D27')
$cell.NumberFormat = "@"
$cell.Value = '9.10'
$cell.ClearFormats()
$ws.Range('E27').Value = '  -1.08%  '
$cell = $ws.Range('D28')
$cell.NumberFormat = "@"
$cell.Value = '582.31'
$cell.ClearFormats()
$ws.Range('E28').Value = '  -8.71%  '
$cell = $ws.Range('D29')
$cell.NumberFormat = "@"
$cell.Value = '2.570.90'
$cell.ClearFormats()
$ws.Range('E29').Value = '  +0.58%  '
$cell = $ws.Range('D30')
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.ClearFormats()
$ws.Range('E30').Value = '  +0.11%  '
$ws.Range('E31').Value = '  -3.40%  '
$cell = $ws.Range('D32')
$cell.NumberFormat = "@"
$cell.Value = '7.92'
$cell.ClearFormats()
$ws.Range('E32').Value = '  -1.96%  '
$ws.Range('E33').Value = '  -5.31%  '
$ws.Range('E34').Value = '  -0.91%  '
$ws.Range('E35').Value = '  -5.72%  '
$ws.Range('E36').Value = '  +0.12%  '
$cell = $ws.Range('D37')
$cell.NumberFormat = "@"
$cell.Value = '4.71'
$cell.ClearFormats()
$ws.Range('E37').Value = '  -5.35%  '
$ws.Range('E38').Value = '  -1.16%  '
$ws.Range('E39').Value = '  -3.83%  '
$cell = $ws.Range('D40')
$cell.NumberFormat = "@"
$cell.Value = '150.92'
$cell.ClearFormats()
$ws.Range('E40').Value = '  +1.25%  '
$cell = $ws.Range('D41')
$cell.NumberFormat = "@"
$cell.Value = '18.27'
$cell.ClearFormats()
$ws.Range('E41').Value = '  -0.89%  '
$ws.Range('E42').Value = '  -2.70%  '
$ws.Range('E44').Value = '  -3.77%  '
$cell = $ws.Range('D45')
$cell.NumberFormat = "@"
$cell.Value = '41.63'
$cell.ClearFormats()
$ws.Range('E45').Value = '  -2.60%  '
$cell = $ws.Range('D46')
$cell.NumberFormat = "@"
$cell.Value = '2.36'
$cell.ClearFormats()
$ws.Range('E46').Value = '  -5.59%  '
$cell = $ws.Range('D47')
$cell.NumberFormat = "@"
$cell.Value = '0.0₆0290'
$cell.ClearFormats()
$ws.Range('E47').Value = '  +24.52%  '
$cell = $ws.Range('D48')
$cell.NumberFormat = "@"
$cell.Value = '142.70'
$cell.ClearFormats()
$ws.Range('E48').Value = '  -1.00%  '
$ws.Range('E49').Value = '  -2.74%  '
$ws.Range('E50').Value = '  -0.05%  '
$cell = $ws.Range('D51')
$cell.NumberFormat = "@"
$cell.Value = '19.65'
$cell.ClearFormats()
$ws.Range('E51').Value = '  -0.67%  '
